# Fellowship Conference 2022 Poster - grammar fix
# "...the brain actively learns how talkers speak, and constructs expectations..."
#   -> "...the brain actively learns how talkers speak, and construct expectations..."
# The sentence is also re-split into three runs at the edit boundaries, matching
# how PowerPoint records an in-place text edit.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Prefer resolving the shape by its stable name; fall back to the known
# top-level shape index if the name lookup isn't available.
try {
    $shp = $s.Shapes.Item("TextBox 59")
} catch {
    $shp = $s.Shapes.Item(10)
}

if ($shp.Name -ne "TextBox 59") {
    throw "Expected to find shape 'TextBox 59' but found '$($shp.Name)'"
}

$tr = $shp.TextFrame.TextRange
$fullText = $tr.Text

$needle = "and constructs expectations"
$idx = $fullText.IndexOf($needle)
if ($idx -lt 0) {
    throw "Could not find target phrase '$needle' in shape text"
}

# Remove the extra "s" in "constructs" -> "construct" (1-based COM index).
$sPos = $idx + "and construct".Length + 1
$sChar = $tr.Characters($sPos, 1)
if ($sChar.Text -ne "s") {
    throw "Unexpected character at computed position: [$($sChar.Text)]"
}
$sChar.Text = ""

# Re-read text/positions after the deletion.
$fullText = $tr.Text
$boundary1 = $idx + 1                       # 1-based start of "and construct "
$boundary2 = $boundary1 + "and construct ".Length   # 1-based start of "expectations..."

# Force the run to split into three runs at the two boundaries, without
# altering any characters, by re-assigning each sub-range's Text to itself.
$runB = $tr.Characters($boundary1, "and construct ".Length)
if ($runB.Text -ne "and construct ") {
    throw "Unexpected run B text: [$($runB.Text)]"
}
$runB.Text = $runB.Text

$runCText = "expectations about how that talker will produce speech in the future. Though this process often occurs without the listener noticing,"
$runC = $tr.Characters($boundary2, $runCText.Length)
if ($runC.Text -ne $runCText) {
    throw "Unexpected run C text: [$($runC.Text)]"
}
$runC.Text = $runC.Text

Write-Output "Final text: $($tr.Text)"
